$d = $word.ActiveDocument

# The target paragraph is the (currently empty) list item that holds the
# "_GoBack" bookmark. Locate it robustly via the bookmark position rather
# than a hard-coded paragraph index.
$bm = $d.Bookmarks.Item("_GoBack")
$bmPos = $bm.Range.Start

$dstIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($bmPos -ge $cand.Range.Start -and $bmPos -lt $cand.Range.End) {
        $dstIndex = $i
    }
}
$pDst = $d.Paragraphs.Item($dstIndex)
$pSrc = $d.Paragraphs.Item($dstIndex - 1)

# Remove the bookmark so it doesn't interfere with the ordering of the
# inserted runs; it is re-added afterwards at its correct position.
$bm.Delete()

$dstRange = $pDst.Range
$dstStart = $dstRange.Start

# Clone the sibling "UpdateDeleteWindow(Delete)" paragraph's run structure
# (3 runs: red "...(", plain word, red ")") so the new runs get identical
# rPr (fonts/size/lang) to the rest of the list.
$dstRange.FormattedText = $pSrc.Range.FormattedText

# Fix up the wording for this list item while keeping the copied run
# formatting intact.
$findRange1 = $pDst.Range
$findRange1.Find.Execute("UpdateDeleteWindow(", $true, $false, $false, $false, $false, $true, 1, $false, "PetUpdateWindow(", 2) | Out-Null

$findRange2 = $pDst.Range
$findRange2.Find.Execute("Delete", $true, $false, $false, $false, $false, $true, 1, $false, "Update", 2) | Out-Null

# Re-insert the "_GoBack" bookmark between "Update" and the closing ")",
# matching its original location in the paragraph.
$bmPos2 = $dstStart + "PetUpdateWindow(".Length + "Update".Length
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos2, $bmPos2)) | Out-Null
